# Update the date on the title slide's subtitle placeholder from
# "March 31, 2016" to "March 29, 2018".
#
# The edit is done as two in-place text replacements (mirroring how a
# user would retype the day and the year directly in PowerPoint), which
# causes PowerPoint to split the single run into three runs:
#   "March " / "29" / ", 2018"

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)          # "Subtitle 2" placeholder
$tr = $sh.TextFrame.TextRange

# "March 31, 2016"
#        ^^ chars 7-8   ^^^^^^ chars 9-14
# Replace "31" -> "29"
$tr.Characters(7, 2).Text = "29"

# Replace ", 2016" -> ", 2018"
$tr.Characters(9, 6).Text = ", 2018"
